$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the duplicated block of columns U:AN (shift remaining cells left).
#    Before the edit, columns U:AN held a third repetition of the
#    "1Pair-A ... MaxUnique" header block (row 2) and the matching header
#    index numbers (row 1). After the edit only columns A:T remain.
$ws.Range("U1:AN19").Delete(-4159)

# 2. Re-order the HKL index labels on row 2 (columns C:J) to the new
#    "Holden" scheme ordering.
$ws.Cells.Item(2,3).Value  = "[3, 2, 1]"
$ws.Cells.Item(2,4).Value  = "[3, 1, 0]"
$ws.Cells.Item(2,5).Value  = "[2, 2, 2]"
$ws.Cells.Item(2,6).Value  = "[1, 1, 0]"
$ws.Cells.Item(2,7).Value  = "[2, 0, 0]"
$ws.Cells.Item(2,8).Value  = "[2, 2, 0]"
$ws.Cells.Item(2,9).Value  = "[4, 0, 0]"
$ws.Cells.Item(2,10).Value = "[2, 1, 1]"

# 3. Append four new simulation rows (20-23) following the same pattern as
#    the existing rows: column A is the zero-based row index (styled like
#    the rest of column A), column B is the scheme name, and columns C:T
#    are all populated with 1.
$newRows = @(
    @(18, "HexGrid-90degTilt2.5degRes"),
    @(19, "HexGrid-90degTilt5degRes"),
    @(20, "HexGrid-90degTilt10degRes"),
    @(21, "HexGrid-90degTilt15degRes")
)

$r = 20
foreach ($pair in $newRows) {
    $idx = $pair[0]
    $name = $pair[1]

    $ws.Cells.Item($r, 1).Value = $idx
    $ws.Cells.Item($r, 2).Value = $name

    for ($c = 3; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }

    $r = $r + 1
}

# Copy the formatting (bold, centered, bordered) used by the rest of
# column A down onto the newly added rows A20:A23.
$ws.Range("A2").Copy()
$ws.Range("A20:A23").PasteSpecial(-4122)
$excel.CutCopyMode = 0
